# Updated cryptos list with GitHub Actions scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.615.87"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "3.075.34"
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("D4").Formula = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Formula = "'589.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Formula = "'155.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.81%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Formula = "'0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").Value = "3.076.20"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").Formula = "'0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.59%  "
$ws.Range("D11").Formula = "'5.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").Formula = "'0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.56%  "
$ws.Range("D13").Formula = "'36.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("D14").Formula = "'0.0000237"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.40%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.582.29"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Formula = "'0.119"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "63.490.15"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Formula = "'7.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.62%  "
$ws.Range("D19").Value = "3.069.54"
$ws.Range("E19").Value = "  -3.01%  "
$ws.Range("D20").Formula = "'472.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Formula = "'14.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").Formula = "'0.705"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.60%  "
$ws.Range("D23").Formula = "'7.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").Formula = "'2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").Formula = "'80.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").Formula = "'12.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").Formula = "'10.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("D28").Formula = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Formula = "'7.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Formula = "'2.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Formula = "'0.996"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").Formula = "'2.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.63%  "
$ws.Range("D33").Formula = "'27.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("D35").Value = "0.0₃0821"
$ws.Range("E35").Value = "  -4.86%  "
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("D37").Formula = "'5.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("D38").Formula = "'3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("D39").Formula = "'2.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.98%  "
$ws.Range("D40").Formula = "'50.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").Formula = "'9.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("D42").Formula = "'436.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.12%  "
$ws.Range("D43").Formula = "'0.286"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("D44").Formula = "'40.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("D45").Formula = "'0.112"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("D46").Formula = "'0.0359"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.60%  "
$ws.Range("D47").Value = "2.796.16"
$ws.Range("E47").Value = "  -4.03%  "
$ws.Range("D48").Formula = "'129.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Formula = "'25.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Formula = "'2.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.86%  "
